# Rename the original sheet to "summary" and add a new "LE" sheet after it.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "summary"

# Re-balance the header alignment on the summary sheet (visual result is
# unchanged, single-column headers stay left-aligned, the merged BD header
# stays centered) -- done via the object model so the style table gets
# rebuilt the same way Excel does it.
$ws1.Range("E1").HorizontalAlignment = -4108
$ws1.Range("H1").HorizontalAlignment = -4108
$ws1.Range("K1").HorizontalAlignment = -4108
$ws1.Range("I1:J1").HorizontalAlignment = -4131

# Add the new "LE" (Local Engagement) sheet right after "summary".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "LE"

# Headers
$ws2.Range("A1").Value = "Region"
$ws2.Range("B1").Value = "Goal"
$ws2.Range("C1").Value = "Status"
$ws2.Range("D1").Value = "Trend"

# Data rows
$ws2.Range("A2").Value = "Hawaii"
$ws2.Range("B2").Value = "LIV"
$ws2.Range("C2").Value = 88.62
$ws2.Range("D2").Value = 0.54

$ws2.Range("A3").Value = "Maui Nui"
$ws2.Range("B3").Value = "LIV"
$ws2.Range("C3").Value = 94.89
$ws2.Range("D3").Value = 0.55

$ws2.Range("A4").Value = "Oahu"
$ws2.Range("B4").Value = "LIV"
$ws2.Range("C4").Value = 75.43
$ws2.Range("D4").Value = 0.51

$ws2.Range("A5").Value = "Kauai"
$ws2.Range("B5").Value = "LIV"
$ws2.Range("C5").Value = 76.11
$ws2.Range("D5").Value = 0.55

# Explicit black font colour on the whole used range (matches the new
# font record added to styles.xml).
$ws2.Range("A1:D5").Font.Color = 0

# Selection on the new sheet.
$ws2.Range("E2").Select()

# Leave "summary" as the active sheet/tab with C3 selected.
$ws1.Activate()
$ws1.Range("C3").Select()
